$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.591.50'
$ws.Range("E2").Value = '  +3.12%  '
# Row 3
$ws.Range("D3").Value = '1.857.91'
$ws.Range("E3").Value = '  +2.02%  '
# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '272.46'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.25%  '
# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5267'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +3.05%  '
# Row 8
$ws.Range("E8").Value = '  -5.04%  '
# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06784'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.61%  '
# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.84'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.07%  '
# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.7921'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -4.52%  '
# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.07733'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.14%  '
# Row 13
$ws.Range("D13").Value = '1.885.02'
$ws.Range("E13").Value = '  +3.51%  '
# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '89.60'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.77%  '
# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.129'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.82%  '
# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.04%  '
# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '14.39'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.72%  '
# Row 18
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.01%  '
# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000007978'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.72%  '
# Row 20
$ws.Range("D20").Value = '26.620.29'
$ws.Range("E20").Value = '  +3.02%  '
# Row 21
$ws.Range("D21").Value = '2.125.19'
$ws.Range("E21").Value = '  +3.71%  '
# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.718'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.84%  '
# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.969'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.45%  '
# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.103'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.27%  '
# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.351'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +4.51%  '
# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '145.64'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +2.44%  '
# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '1.656'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.90%  '
# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '17.18'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.14%  '
# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '111.98'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.11%  '
# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.296'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.88%  '
# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.288'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.90%  '
# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.08891'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.45%  '
# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.04907'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.13%  '
# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.155'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.29%  '
# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.7258'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.97%  '
# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.882'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.22%  '
# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.221'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.49%  '
# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.314'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.74%  '
# Row 39
$ws.Range("E39").Value = '  -1.09%  '
# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.5073'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.15%  '
# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.9360'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.32%  '
# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '115.86'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.64%  '
# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '6.129'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.03%  '
# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '7.985'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.84%  '
# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.07%  '
# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4401'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -3.55%  '
# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.1322'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.65%  '
# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '9.270'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.75%  '
# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '35.96'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.90%  '
# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.05928'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.53%  '
# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.470'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.98%  '

Write-Output "Applied all changes"